$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "057-82511516"
$ws.Range("B2").Value = 2120668841
$ws.Range("C2").Value = 2120668841

# Row 3
$ws.Range("A3").Value = "057-19009594"
$ws.Range("B3").Value = 2061841172
$ws.Range("C3").Value = 2061841172

# Row 4
$ws.Range("A4").Value = "057-83699442"
$ws.Range("B4").Value = 2550335038
$ws.Range("C4").Value = 2550335038

# Row 5
$ws.Range("A5").Value = "057-83699453"
$ws.Range("B5").Value = 2550334549
$ws.Range("C5").Value = 2550334549

# Row 6
$ws.Range("A6").Value = "057-83699453"
$ws.Range("B6").Value = 2550334546
$ws.Range("C6").Value = 2550334546

# Row 7
$ws.Range("A7").Value = "057-83699442"
$ws.Range("B7").Value = 2550335047
$ws.Range("C7").Value = 2550335047

# Row 8 (new)
$ws.Range("A8").Value = "057-82511516"
$ws.Range("B8").Value = 2120668847
$ws.Range("C8").Value = 2120668847

# Row 9 (new, formerly row 5)
$ws.Range("A9").Value = "057-96206751"
$ws.Range("B9").Value = "24S0041891"
$ws.Range("C9").Value = "24S0041891"

# Row 10 (new, formerly row 6)
$ws.Range("A10").Value = "057-83419232"
$ws.Range("B10").Value = "21A0406836"
$ws.Range("C10").Value = "21A0406836"

# Row 11 (new, formerly row 7)
$ws.Range("A11").Value = "057-70312594"
$ws.Range("B11").Value = "DJMKEA4230073"
$ws.Range("C11").Value = 222323293

# Update the selection to reflect the new data extent
[void]$ws.Range("A2:C11").Select()
